$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 17. Existing rows 17-30 shift down to 18-31.
$ws.Rows("17:17").Insert()

# Populate the newly inserted row 17 with the weekly price-report entry.
# (Mercado/Region/Categoria/Variedad/Calidad/Unidad/Origen/Kg-Unidades/Clasificacion
# columns carry the same values as the rest of this market's rows; only the
# date and price/volume figures differ for this new weekly record.)
$ws.Range("A17").Value = 1
$ws.Range("B17").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C17").Value = "Arica y Parinacota"
$ws.Range("D17").Value = 44827
$ws.Range("E17").Value = 15
$ws.Range("F17").Value = 100112013
$ws.Range("G17").Value = "Alcachofa"
$ws.Range("H17").Value = "Madrigal"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 100
$ws.Range("K17").Value = 14000
$ws.Range("L17").Value = 15000
$ws.Range("M17").Value = 14500
$ws.Range("N17").Value = "$/caja 40 unidades"
$ws.Range("O17").Value = "Región de Coquimbo"
$ws.Range("P17").Value = 362
$ws.Range("Q17").Value = 40
$ws.Range("R17").Value = "Hortaliza"
